# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.367.34"
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = "'2.427.34"
$ws.Range("E3").Value = '  +1.17%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'563.37"
$ws.Range("E5").Value = '  +2.37%  '
$ws.Range("D6").Value = "'166.67"
$ws.Range("E6").Value = '  +5.59%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.96%  '
$ws.Range("E9").Value = '  +8.05%  '
$ws.Range("D10").Value = "'2.425.52"
$ws.Range("E10").Value = '  +1.26%  '
$ws.Range("E11").Value = '  -1.69%  '
$ws.Range("E12").Value = '  +2.26%  '
$ws.Range("E13").Value = '  -0.87%  '
$ws.Range("E14").Value = '  +6.20%  '
$ws.Range("D15").Value = "'69.259.88"
$ws.Range("E15").Value = '  +2.24%  '
$ws.Range("D16").Value = "'2.873.75"
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").Value = "'23.91"
$ws.Range("E17").Value = '  +4.99%  '
$ws.Range("D18").Value = "'2.435.35"
$ws.Range("E18").Value = '  +2.26%  '
$ws.Range("D19").Value = "'10.81"
$ws.Range("E19").Value = '  +5.28%  '
$ws.Range("D20").Value = "'342.98"
$ws.Range("E20").Value = '  +4.38%  '
$ws.Range("E21").Value = '  +6.12%  '
$ws.Range("D22").Value = "'3.86"
$ws.Range("E22").Value = '  +3.01%  '
$ws.Range("E23").Value = '  +7.09%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").Value = "'66.02"
$ws.Range("E25").Value = '  +0.65%  '
$ws.Range("E26").Value = '  +6.16%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").Value = "'1.04"
$ws.Range("E27").Value = '  +4.22%  '
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value = "'8.48"
$ws.Range("E28").Value = '  +6.45%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = "'2.552.84"
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("E30").Value = '  +7.38%  '
$ws.Range("E31").Value = '  +5.74%  '
$ws.Range("E32").Value = '  +11.13%  '
$ws.Range("D33").Value = "'453.84"
$ws.Range("E33").Value = '  +9.18%  '
$ws.Range("E34").Value = '  +0.07%  '
$ws.Range("D35").Value = "'1.62"
$ws.Range("E35").Value = '  +2.94%  '
$ws.Range("D36").Value = "'159.10"
$ws.Range("E36").Value = '  +1.35%  '
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("E38").Value = '  +7.10%  '
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").Value = "'18.24"
$ws.Range("E40").Value = '  +3.73%  '
$ws.Range("E41").Value = '  +4.35%  '
$ws.Range("E42").Value = '  +4.47%  '
$ws.Range("E43").Value = '  +5.13%  '
$ws.Range("D44").Value = "'37.82"
$ws.Range("E44").Value = '  +1.92%  '
$ws.Range("E45").Value = '  +3.34%  '
$ws.Range("E46").Value = '  +8.32%  '
$ws.Range("D47").Value = "'135.80"
$ws.Range("E47").Value = '  +6.54%  '
$ws.Range("E48").Value = '  +3.96%  '
$ws.Range("E49").Value = '  +2.91%  '
$ws.Range("E50").Value = '  +3.89%  '
$ws.Range("E51").Value = '  +2.76%  '
